$wb = $excel.ActiveWorkbook

# Update "展览" sheet (sheetId 1)
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 4829
$wsExhibition.Range("F3").Value = 143
$wsExhibition.Range("F4").Value = 844

# Update "全部类型" sheet (sheetId 4)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 4829
$wsAll.Range("F3").Value = 143
$wsAll.Range("F4").Value = 844
